# Change values in database
# Updates production figures (column I) on the "1 нф" sheet and moves the
# active cell selection, matching the authoring edit captured in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1 нф")
$ws.Activate()

$ws.Range("I3").Value = 1200
$ws.Range("I4").Value = 1270
$ws.Range("I5").Value = 1176

$ws.Range("I7").Value = 12503
$ws.Range("I8").Value = 11300
$ws.Range("I9").Value = 11670

$ws.Range("I10").Select()
